$wb = $excel.ActiveWorkbook

# --- Sheet 1: TextBoxData ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Columns.Item(1).Insert()
$ws1.Range("A1").Value = "NameId"
for ($i = 0; $i -lt 20; $i++) {
    $row = 2 + $i
    $ws1.Cells.Item($row, 1).Value = 101011 + $i
}

# --- Sheet 2: PracticeFormData ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Columns.Item(1).Insert()
$ws2.Range("A1").Value = "FormId"
for ($i = 0; $i -lt 20; $i++) {
    $row = 2 + $i
    $ws2.Cells.Item($row, 1).Value = 101011 + $i
}

# Fill the new last column (N) UploadFilePath for rows 7-21 with the same
# path already present in row 6 (column N after the insert).
$lastPath = $ws2.Cells.Item(6, 14).Value()
for ($row = 7; $row -le 21; $row++) {
    $ws2.Cells.Item($row, 14).Value = $lastPath
}

# --- Selections / active sheet ---
$ws1.Range("A2").Select()
$ws2.Activate()
$ws2.Range("M4").Select()

Write-Host "done"
